$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 23187.166
$ws.Range("I40").Value = 6464
$ws.Range("J40").Value = 46599.6
$ws.Range("K40").Value = 6464
$ws.Range("L40").Value = 46599.6
$ws.Range("M40").Value = -6289
$ws.Range("N40").Value = -46949.6
$ws.Range("H43").Value = 14666.546
$ws.Range("I43").Value = 19936.4
$ws.Range("J43").Value = 10275
$ws.Range("K43").Value = 19936.4
$ws.Range("L43").Value = 10275
$ws.Range("M43").Value = -19867.4
$ws.Range("N43").Value = -10413
$ws.Range("H51").Value = 23608.908
$ws.Range("I51").Value = 9999.5
$ws.Range("J51").Value = 26633.223
$ws.Range("K51").Value = 9999.5
$ws.Range("L51").Value = 26633.223
$ws.Range("M51").Value = -9515.5
$ws.Range("N51").Value = -27601.223
$ws.Range("H94").Value = 645.5714
$ws.Range("I94").Value = 645.5714
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 645.5714
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -194.5714
$ws.Range("N94").ClearContents()
$ws.Range("H97").Value = 3527
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H100").Value = 3384.2856
$ws.Range("I100").Value = 3735.25
$ws.Range("J100").Value = 2916.3333
$ws.Range("K100").Value = 3735.25
$ws.Range("L100").Value = 2916.3333
$ws.Range("M100").Value = -3194.25
$ws.Range("N100").Value = -3998.3333
$ws.Range("H112").Value = 3814.5557
$ws.Range("I112").Value = 1366.6666
$ws.Range("J112").Value = 5038.5
$ws.Range("K112").Value = 4099.9998
$ws.Range("L112").Value = 15115.5
$ws.Range("M112").Value = -2991.9998
$ws.Range("N112").Value = -17331.5
$ws.Range("H127").Value = 674.25
$ws.Range("I127").Value = 674.25
$ws.Range("K127").Value = 2022.75
$ws.Range("M127").Value = 2937.25
$ws.Range("H129").Value = 1631.1666
$ws.Range("I129").Value = 1631.1666
$ws.Range("K129").Value = 4893.4998
$ws.Range("M129").Value = 106.5002000000004
$ws.Range("H131").Value = 836.6
$ws.Range("I131").Value = 836.6
$ws.Range("K131").Value = 2509.8
$ws.Range("M131").Value = 2530.2
$ws.Range("H132").Value = 3831.8096
$ws.Range("I132").Value = 4014.7222
$ws.Range("J132").Value = 2734.3333
$ws.Range("K132").Value = 12044.1666
$ws.Range("L132").Value = 8202.999899999999
$ws.Range("M132").Value = -9514.1666
$ws.Range("N132").Value = -13262.9999
$ws.Range("H141").Value = 5210
$ws.Range("I141").Value = 4781.6665
$ws.Range("K141").Value = 14344.9995
$ws.Range("M141").Value = -9164.999500000002
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1907.1666
$ws.Range("I2").Value = 2138.6
$ws.Range("J2").Value = 750
$ws.Range("K2").Value = 2138.6
$ws.Range("L2").Value = 750
$ws.Range("M2").Value = -2025.6
$ws.Range("N2").Value = -976
$ws.Range("H45").Value = 8516.5
$ws.Range("I45").Value = 12782.211
$ws.Range("K45").Value = 12782.211
$ws.Range("M45").Value = -12405.211
$ws.Range("H61").Value = 7474.269
$ws.Range("I61").Value = 7470.5557
$ws.Range("J61").Value = 7482.625
$ws.Range("K61").Value = 7470.5557
$ws.Range("L61").Value = 7482.625
$ws.Range("M61").Value = -7258.5557
$ws.Range("N61").Value = -7906.625
$ws.Range("H63").Value = 1858.1
$ws.Range("I63").Value = 953.44446
$ws.Range("K63").Value = 953.44446
$ws.Range("M63").Value = -267.44446
$ws.Range("H66").Value = 1858.1
$ws.Range("I66").Value = 953.44446
$ws.Range("K66").Value = 4767.2223
$ws.Range("M66").Value = -1335.2223
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H116").Value = 1907.1666
$ws.Range("I116").Value = 2138.6
$ws.Range("J116").Value = 750
$ws.Range("K116").Value = 2138.6
$ws.Range("L116").Value = 750
$ws.Range("M116").Value = 155.4000000000001
$ws.Range("N116").Value = -5338
$ws.Range("H132").Value = 3227.913
$ws.Range("I132").Value = 2986.55
$ws.Range("K132").Value = 8959.650000000001
$ws.Range("M132").Value = -6429.650000000001
$ws.Range("H136").Value = 7474.269
$ws.Range("I136").Value = 7470.5557
$ws.Range("J136").Value = 7482.625
$ws.Range("K136").Value = 22411.6671
$ws.Range("L136").Value = 22447.875
$ws.Range("M136").Value = -19861.6671
$ws.Range("N136").Value = -27547.875
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1907.1666
$ws.Range("I3").Value = 2138.6
$ws.Range("J3").Value = 750
$ws.Range("K3").Value = 2138.6
$ws.Range("L3").Value = 750
$ws.Range("M3").Value = -2024.6
$ws.Range("N3").Value = -978
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4329.0625
$ws.Range("I31").Value = 4532.5
$ws.Range("J31").Value = 3718.75
$ws.Range("K31").Value = 4532.5
$ws.Range("L31").Value = 3718.75
$ws.Range("M31").Value = -4237.5
$ws.Range("N31").Value = -4308.75
$ws.Range("H34").Value = 4329.0625
$ws.Range("I34").Value = 4532.5
$ws.Range("J34").Value = 3718.75
$ws.Range("K34").Value = 4532.5
$ws.Range("L34").Value = 3718.75
$ws.Range("M34").Value = -4330.5
$ws.Range("N34").Value = -4122.75
$ws.Range("H59").Value = 92856.42999999999
$ws.Range("J59").Value = 99999.164
$ws.Range("L59").Value = 99999.164
$ws.Range("N59").Value = -102289.164
$ws.Range("H62").Value = 2862
$ws.Range("I62").Value = 2724.5
$ws.Range("K62").Value = 2724.5
$ws.Range("M62").Value = -2100.5
$ws.Range("H65").Value = 2862
$ws.Range("I65").Value = 2724.5
$ws.Range("K65").Value = 13622.5
$ws.Range("M65").Value = -10502.5
$ws.Range("H99").Value = 5024
$ws.Range("I99").Value = 3366
$ws.Range("K99").Value = 3366
$ws.Range("M99").Value = -1868
$ws.Range("H126").Value = 5024
$ws.Range("I126").Value = 3366
$ws.Range("K126").Value = 10098
$ws.Range("M126").Value = -7628
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1523.8182
$ws.Range("J129").Value = 1996.8
$ws.Range("L129").Value = 5990.4
$ws.Range("N129").Value = -15990.4
$ws.Range("H131").Value = 2414.652
$ws.Range("I131").Value = 1254
$ws.Range("J131").Value = 3307.4614
$ws.Range("K131").Value = 3762
$ws.Range("L131").Value = 9922.3842
$ws.Range("M131").Value = 1278
$ws.Range("N131").Value = -20002.3842
$ws.Range("H139").Value = 7013.476
$ws.Range("I139").Value = 3814.111
$ws.Range("J139").Value = 9413
$ws.Range("K139").Value = 11442.333
$ws.Range("L139").Value = 28239
$ws.Range("M139").Value = -6302.332999999999
$ws.Range("N139").Value = -38519
$ws.Range("H140").Value = 2503.3809
$ws.Range("I140").Value = 2046.0625
$ws.Range("J140").Value = 3966.8
$ws.Range("K140").Value = 6138.1875
$ws.Range("L140").Value = 11900.4
$ws.Range("M140").Value = -958.1875
$ws.Range("N140").Value = -22260.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 723.2308
$ws.Range("I97").Value = 695.1429000000001
$ws.Range("K97").Value = 695.1429000000001
$ws.Range("M97").Value = -199.1429000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 10475.25
$ws.Range("I2").Value = 10299.667
$ws.Range("K2").Value = 10299.667
$ws.Range("M2").Value = -10187.667
$ws.Range("H7").Value = 15507.315
$ws.Range("J7").Value = 8823.5
$ws.Range("L7").Value = 8823.5
$ws.Range("N7").Value = -9047.5
$ws.Range("H61").Value = 2345.8333
$ws.Range("I61").Value = 2635.625
$ws.Range("J61").Value = 1766.25
$ws.Range("K61").Value = 2635.625
$ws.Range("L61").Value = 1766.25
$ws.Range("M61").Value = -2433.625
$ws.Range("N61").Value = -2170.25
$ws.Range("H113").Value = 2345.8333
$ws.Range("I113").Value = 2635.625
$ws.Range("J113").Value = 1766.25
$ws.Range("K113").Value = 2635.625
$ws.Range("L113").Value = 1766.25
$ws.Range("M113").Value = -465.625
$ws.Range("N113").Value = -6106.25
$ws.Range("H126").Value = 15507.315
$ws.Range("J126").Value = 8823.5
$ws.Range("L126").Value = 26470.5
$ws.Range("N126").Value = -31410.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2203.5386
$ws.Range("I100").Value = 2176.6667
$ws.Range("K100").Value = 4353.3334
$ws.Range("M100").Value = -3812.3334
$ws.Range("H119").Value = 6000000
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H123").Value = 76164.5
$ws.Range("I123").Value = 71900
$ws.Range("J123").Value = 80429
$ws.Range("K123").Value = 71900
$ws.Range("L123").Value = 80429
$ws.Range("M123").Value = -67000
$ws.Range("N123").Value = -90229
$ws.Range("H132").Value = 3387.1853
$ws.Range("I132").Value = 3602.8
$ws.Range("J132").Value = 2771.1428
$ws.Range("K132").Value = 10808.4
$ws.Range("L132").Value = 8313.428400000001
$ws.Range("M132").Value = -8278.400000000001
$ws.Range("N132").Value = -13373.4284
$ws.Range("H136").Value = 2352.973
$ws.Range("I136").Value = 1817.6666
$ws.Range("K136").Value = 5452.9998
$ws.Range("M136").Value = -2902.9998
